$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the tag text for "Reverse Linked List" (row 97, column C)
$ws.Range("C97").Value = "#linked-list #recursive #必背"

# 2. Append a new problem row for "239. Sliding Window Maximum", copying the
#    formatting (styles/number formats) from the row above it first so no new
#    cell styles are introduced.
$ws.Range("A97:I97").Copy()
$ws.Range("A99:I99").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A99").Value = 239
$ws.Range("B99").Value = "Sliding Window Maximum"
$ws.Range("C99").Value = "#sliding-window #monotonic-queue #array #queue #核心 "
$ws.Range("D99").Value = "hard"
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 2
$ws.Range("G99").Value = 40
$ws.Range("H99").Value = 45864
$ws.Range("I99").Value = 45864

$ws.Rows.Item(99).RowHeight = 68

# Update the sheet view so the newly added row is visible/selected
$ws.Range("H99:I99").Select()
